$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Flat, Flags<0>" --- (added first so "Order int" gets the lower shared-string index)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A7").Value = "Order int"
$ws2.Range("D7").Value = 332
$ws2.Range("E7").Formula = "=(D7/D`$2)-1"
$ws2.Range("E7").Style = $ws2.Range("E6").Style
$ws2.Range("E7").NumberFormat = $ws2.Range("E6").NumberFormat
$ws2.Range("F7").Formula = "=(D7/D6)-1"
$ws2.Range("F7").Style = $ws2.Range("F6").Style
$ws2.Range("F7").NumberFormat = $ws2.Range("F6").NumberFormat
$ws2.Range("A7").Select()

# --- Sheet 1: "Tex, Flags<0>" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A13").Value = "Order"
$ws1.Range("D13").Value = 98
$ws1.Range("E13").Formula = "=(D13/D`$2)-1"
$ws1.Range("E13").Style = $ws1.Range("E6").Style
$ws1.Range("E13").NumberFormat = $ws1.Range("E6").NumberFormat
$ws1.Range("E13").Select()

# --- Sheet 3: "V1" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D6").Value = 332
$ws3.Range("D2").Select()
